$d = $word.ActiveDocument

# 1. Merge "Response " + "Code/Message" into a single run (text unchanged)
$d.Content.Find.Execute("Response Code/Message", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Response Code/Message", 2)

# 2. Merge "stored with /t2" + "_update" + " endpoint" into a single run (text unchanged)
$d.Content.Find.Execute("stored with /t2_update endpoint", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 0)

# 3. Change heading "LED Parameters/Arguments" to "Custom Parameters/Arguments"
$d.Content.Find.Execute("LED Parameters/Arguments", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Custom Parameters/Arguments", 2)

# 4. Remove the bookmark near the end of the document
foreach ($bm in $d.Bookmarks) {
    if ($bm.Name -eq "_GoBack") {
        $bm.Delete()
    }
}

# 5. Re-add the bookmark right after "Custom" in the heading
$rng = $d.Content
$rng.Find.Execute("Custom Parameters/Arguments", $true, $false, $false, $false, $false,
                   $true, 1, $false, "", 0)
if ($rng.Find.Found) {
    $bmRange = $d.Range($rng.Start, $rng.Start + 6)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
